# Add a "fictive_time" column to the Results schema, inserted right before
# the existing "hardware" column (which, together with memory_in_kb and
# time_in_s, shifts one column to the right).
$wb = $excel.ActiveWorkbook

foreach ($name in @("Results", "Results1")) {
    $ws = $wb.Sheets.Item($name)
    # xlShiftToRight = -4161
    $ws.Columns.Item(4).Insert(-4161)
    $ws.Range("D1").Value = "fictive_time"
}
